$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.811.40'
$ws.Range("E2").Value = '  -5.00%  '
$ws.Range("D3").Value = '2.991.06'
$ws.Range("E3").Value = '  -5.18%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '542.04'
$ws.Range("E5").Value = '  -5.77%  '
$ws.Range("D6").Value = '151.59'
$ws.Range("E6").Value = '  -9.36%  '
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = '0.571'
$ws.Range("E8").Value = '  -1.27%  '
$ws.Range("D9").Value = '3.002.27'
$ws.Range("E9").Value = '  -5.41%  '
$ws.Range("E10").Value = '  -4.97%  '
$ws.Range("D11").Value = '6.13'
$ws.Range("E11").Value = '  -8.00%  '
$ws.Range("E12").Value = '  -4.22%  '
$ws.Range("D13").Value = '3.512.09'
$ws.Range("E13").Value = '  -5.20%  '
$ws.Range("E14").Value = '  -1.19%  '
$ws.Range("D15").Value = '61.829.05'
$ws.Range("E15").Value = '  -4.97%  '
$ws.Range("D16").Value = '23.99'
$ws.Range("E16").Value = '  -5.07%  '
$ws.Range("D17").Value = '2.998.11'
$ws.Range("E17").Value = '  -4.70%  '
$ws.Range("E18").Value = '  -6.33%  '
$ws.Range("D19").Value = '5.17'
$ws.Range("E19").Value = '  -2.00%  '
$ws.Range("E20").Value = '  -4.54%  '
$ws.Range("D21").Value = '377.68'
$ws.Range("E21").Value = '  -8.98%  '
$ws.Range("D22").Value = '6.69'
$ws.Range("E22").Value = '  -5.79%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("D24").Value = '66.19'
$ws.Range("E24").Value = '  -4.19%  '
$ws.Range("D25").Value = '3.114.78'
$ws.Range("E25").Value = '  -4.84%  '
$ws.Range("E26").Value = '  -3.86%  '
$ws.Range("E27").Value = '  -4.13%  '
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  +0.37%  '
$ws.Range("D29").Value = '0.0₃0934'
$ws.Range("E29").Value = '  -10.95%  '
$ws.Range("D30").Value = '8.24'
$ws.Range("E30").Value = '  -10.92%  '
$ws.Range("D31").Value = '0.999'
$ws.Range("E31").Value = '  +0.00%  '
$ws.Range("E32").Value = '  -5.30%  '
$ws.Range("D33").Value = '20.47'
$ws.Range("E33").Value = '  -4.76%  '
$ws.Range("D34").Value = '160.64'
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("E35").Value = '  -5.37%  '
$ws.Range("D36").Value = '4.60'
$ws.Range("E36").Value = '  -7.97%  '
$ws.Range("E37").Value = '  -6.86%  '
$ws.Range("E38").Value = '  -7.37%  '
$ws.Range("E39").Value = '  -8.65%  '
$ws.Range("D40").Value = '37.61'
$ws.Range("E40").Value = '  -2.18%  '
$ws.Range("D41").Value = '2.415.54'
$ws.Range("E41").Value = '  -7.78%  '
$ws.Range("E42").Value = '  -6.72%  '
$ws.Range("D43").Value = '22.06'
$ws.Range("E43").Value = '  -8.50%  '
$ws.Range("E44").Value = '  -3.74%  '
$ws.Range("D45").Value = '0.0590'
$ws.Range("E45").Value = '  -5.44%  '
$ws.Range("E46").Value = '  -3.64%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  -5.53%  '
$ws.Range("E49").Value = '  -2.81%  '
$ws.Range("B50").Value = 'Bittensor'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D50").Value = '266.57'
$ws.Range("E50").Value = '  -8.90%  '
$ws.Range("B51").Value = 'InjectiveProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D51").Value = '19.63'
$ws.Range("E51").Value = '  -9.36%  '
